$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (dSF) value updates per row, per the repull/mean-calculation fix
$updates = @{
    2  = -5
    4  = -4
    5  = -7
    6  = -5
    7  = -6
    8  = 2
    10 = -2
    11 = 2
    12 = -1
    14 = 4
    16 = 2
    17 = -4
    19 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
